{"js": "// Replace the arithmetic-problem text in every cell of the worksheet's\n// (single) table with the new set of problems, cell-for-cell, preserving\n// the existing table/row/cell/paragraph/run formatting.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst newValues = [\n  [\"67-43=\", \"57+3=\", \"41+44=\", \"71-29=\", \"96-16=\"],\n  [\"34+26=\", \"19+3=\", \"12+21=\", \"14+46=\", \"96-12=\"],\n  [\"75-31=\", \"88-76=\", \"44+53=\", \"40+50=\", \"44+1=\"],\n  [\"35+3=\", \"89-30=\", \"73-72=\", \"86-50=\", \"69-6=\"],\n  [\"24+54=\", \"94-85=\", \"69-6=\", \"37+8=\", \"11+21=\"],\n  [\"59+32=\", \"35-18=\", \"46+33=\", \"73-59=\", \"7+87=\"],\n  [\"60+28=\", \"92-68=\", \"20+9=\", \"68+7=\", \"84-55=\"],\n  [\"16+67=\", \"80-23=\", \"0+43=\", \"73-5=\", \"28-17=\"],\n  [\"48+10=\", \"57+31=\", \"14+60=\", \"57-38=\", \"70-68=\"],\n  [\"55-2=\", \"75+8=\", \"61-41=\", \"23+28=\", \"84-73=\"],\n  [\"88-25=\", \"58-44=\", \"94-88=\", \"56-22=\", \"36+23=\"],\n  [\"55+27=\", \"50+39=\", \"77-5=\", \"86-30=\", \"53-35=\"],\n  [\"29-14=\", \"75+23=\", \"79-66=\", \"73-15=\", \"33+54=\"],\n  [\"45+45=\", \"69+18=\", \"36+46=\", \"83-62=\", \"59-24=\"],\n  [\"27+40=\", \"44-12=\", \"56+13=\", \"56-37=\", \"90-77=\"],\n  [\"9+32=\", \"36-30=\", \"82-36=\", \"91-88=\", \"23-3=\"],\n  [\"81-60=\", \"56+25=\", \"82-12=\", \"67-64=\", \"52+22=\"],\n  [\"47-35=\", \"39+1=\", \"15+62=\", \"6+60=\", \"40+45=\"],\n  [\"68-5=\", \"27+12=\", \"60-16=\", \"31+20=\", \"96-87=\"],\n  [\"88-29=\", \"22+37=\", \"18+46=\", \"30+25=\", \"91-91=\"],\n];\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# Replace the arithmetic-problem text in every cell of the worksheet's\n# (single) table with the new set of problems, cell-for-cell, preserving\n# the existing table/row/cell/paragraph/run formatting.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$newValues = @(\n  @(\"67-43=\", \"57+3=\", \"41+44=\", \"71-29=\", \"96-16=\"),\n  @(\"34+26=\", \"19+3=\", \"12+21=\", \"14+46=\", \"96-12=\"),\n  @(\"75-31=\", \"88-76=\", \"44+53=\", \"40+50=\", \"44+1=\"),\n  @(\"35+3=\", \"89-30=\", \"73-72=\", \"86-50=\", \"69-6=\"),\n  @(\"24+54=\", \"94-85=\", \"69-6=\", \"37+8=\", \"11+21=\"),\n  @(\"59+32=\", \"35-18=\", \"46+33=\", \"73-59=\", \"7+87=\"),\n  @(\"60+28=\", \"92-68=\", \"20+9=\", \"68+7=\", \"84-55=\"),\n  @(\"16+67=\", \"80-23=\", \"0+43=\", \"73-5=\", \"28-17=\"),\n  @(\"48+10=\", \"57+31=\", \"14+60=\", \"57-38=\", \"70-68=\"),\n  @(\"55-2=\", \"75+8=\", \"61-41=\", \"23+28=\", \"84-73=\"),\n  @(\"88-25=\", \"58-44=\", \"94-88=\", \"56-22=\", \"36+23=\"),\n  @(\"55+27=\", \"50+39=\", \"77-5=\", \"86-30=\", \"53-35=\"),\n  @(\"29-14=\", \"75+23=\", \"79-66=\", \"73-15=\", \"33+54=\"),\n  @(\"45+45=\", \"69+18=\", \"36+46=\", \"83-62=\", \"59-24=\"),\n  @(\"27+40=\", \"44-12=\", \"56+13=\", \"56-37=\", \"90-77=\"),\n  @(\"9+32=\", \"36-30=\", \"82-36=\", \"91-88=\", \"23-3=\"),\n  @(\"81-60=\", \"56+25=\", \"82-12=\", \"67-64=\", \"52+22=\"),\n  @(\"47-35=\", \"39+1=\", \"15+62=\", \"6+60=\", \"40+45=\"),\n  @(\"68-5=\", \"27+12=\", \"60-16=\", \"31+20=\", \"96-87=\"),\n  @(\"88-29=\", \"22+37=\", \"18+46=\", \"30+25=\", \"91-91=\")\n)\n\nfor ($r = 1; $r -le $newValues.Length; $r++) {\n    $row = $newValues[$r - 1]\n    for ($c = 1; $c -le $row.Length; $c++) {\n        $t.Cell($r, $c).Range.Text = $row[$c - 1]\n    }\n}\n"}
